# Bugfixed the naive forecaster component module
#
# Column A held Excel date serials (with a bespoke "YYYY-MM-DD HH:MM:SS"
# number format) for what are really just annual Q4 observations. Replace
# them with plain "<year>Q4" text labels so the naive forecaster reads the
# period as a label instead of (mis-)parsing it as a timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds 1987Q4, and each subsequent row advances one year through
# row 39 (2024Q4) - mirrors the original A2:A39 date-serial sequence
# (one entry per year, Q4/year-end).
$firstYear = 1987
$lastRow = 39

for ($row = 2; $row -le $lastRow; $row++) {
    $year = $firstYear + ($row - 2)
    $ws.Range("A$row").Value = "${year}Q4"
}

# The old per-cell style (index 2) only existed to carry the custom date
# number format. Now that column A is plain text, reuse the header's style
# (A1) for A2:A39 so the bespoke date format - and the xf that applied it -
# become dead weight instead of staying wired to live cells.
$ws.Range("A1").Copy()
$ws.Range("A2:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
